# Update the "Annual Population Survey" employment rows (Employment volumes,
# Employment by occupation, Employment by industry) with the latest nomis
# release periods: the data has moved on one quarter and the permalink text
# shown in column D has been refreshed to the newer ILR publication date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C4").Value = "Oct 2022 - Sep 2023 (16/01/24)"
$ws.Range("D2:D4").Value = "Jan 2023 - Dec 2023 (16/04/24)"

# Leave the selection where the analyst was last working, matching the
# refreshed data table rows.
$ws.Range("D2:D4").Select()
